# Actualización automática 2025-09-11 12:20:10
#
# Inserts a new advisor/client row ("COBO FOLLECO JORGE ERNESTO") right
# after "CARRION ALVAREZ MARIO ANDRES" (row 15) in both the
# "VENTAS POR GRUPO" and "VENTA MENSUAL" sheets, pushing the existing
# rows 15-40 down to 16-41 and the trailing totals row from 41 to 42.
# The totals row in "VENTAS POR GRUPO" also has its "X de 39" counter
# labels refreshed to "X de 40" to reflect the extra data row.

$wb = $excel.ActiveWorkbook

# ---- Sheet 1: "VENTAS POR GRUPO" (columns A:R) ----------------------------
$ws1 = $wb.Worksheets.Item(1)

$ws1.Rows.Item(15).Insert()
$ws1.Range("A15").Value = "OFICINA-CATAECSA"
$ws1.Range("B15").Value = "COBO FOLLECO JORGE ERNESTO"
for ($col = 3; $col -le 18; $col++) {
    $ws1.Cells.Item(15, $col).Value = 0
}

# Refresh the "X de 39" -> "X de 40" summary labels on the (now shifted)
# totals row 42.
$ws1.Range("C42").Value = "0 de 40"
$ws1.Range("D42").Value = "1 de 40"
$ws1.Range("E42").Value = "1 de 40"
$ws1.Range("F42").Value = "0 de 40"
$ws1.Range("G42").Value = "0 de 40"
$ws1.Range("H42").Value = "0 de 40"
$ws1.Range("I42").Value = "0 de 40"
$ws1.Range("J42").Value = "0 de 40"
$ws1.Range("K42").Value = "0 de 40"
$ws1.Range("L42").Value = "3 de 40"
$ws1.Range("M42").Value = "3 de 40"
$ws1.Range("N42").Value = "0 de 40"
$ws1.Range("O42").Value = "0 de 40"
$ws1.Range("P42").Value = "0 de 40"
$ws1.Range("Q42").Value = "0 de 40"
$ws1.Range("R42").Value = "0 de 40"

# ---- Sheet 2: "VENTA MENSUAL" (columns A:G) --------------------------------
$ws2 = $wb.Worksheets.Item(2)

$ws2.Rows.Item(15).Insert()
$ws2.Range("A15").Value = "OFICINA-CATAECSA"
$ws2.Range("B15").Value = "COBO FOLLECO JORGE ERNESTO"
for ($col = 3; $col -le 7; $col++) {
    $ws2.Cells.Item(15, $col).Value = 0
}
